# Add a new reference row (Jansen, Zuidema, Anten & Martinez-Ramos 2012 -
# Chamaedorea elegans) to the "Original reference list" sheet, immediately
# above the existing "Loreau ... J Ecol 2011" entry (old row 37), pushing
# every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Original reference list")
$ws.Activate()

# Insert a brand-new row at position 37 (this shifts rows 37-97 down to 38-98)
$ws.Rows.Item(37).Insert()

# Fill in the bibliographic details for the new entry
$ws.Cells.Item(37, 1).Value = "Jansen, Zuidema, Anten & Martinez-Ramos"
$ws.Cells.Item(37, 2).Value = "J Ecol"
$ws.Cells.Item(37, 3).Value = 2012
$ws.Cells.Item(37, 4).Value = "Plantae"
$ws.Cells.Item(37, 5).Value = "Chamaedorea elegans"
$ws.Cells.Item(37, 5).Font.Italic = $true
$ws.Cells.Item(37, 7).Value = "Jansen M, Zuidema PA, Anten NPR & Martinez-Ramos M (2012) Strong persistent growth differences govern individual performance and populaton dynamics in a tropical forest understory palm. Journal of Ecoloy 100: 1224-1232"

# Restore the usual selection on this sheet (the formerly-selected cell G59
# is now one row further down, at G38, after the insertion)
$ws.Range("G38").Select()
